$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are treated as plain text so that values
# like "44.185.25" or "87.50" are preserved exactly (no numeric coercion).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.185.25"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.258.14"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.37"
$ws.Range("E5").Value = "  -4.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.23"
$ws.Range("E6").Value = "  -2.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.576"
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  -3.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.81"
$ws.Range("E10").Value = "  -4.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0825"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.38"
$ws.Range("E12").Value = "  -4.58%  "
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.600.11"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.844"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.254.29"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.93"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.066.17"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.98"
$ws.Range("E19").Value = "  -6.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0978"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.38"
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.62"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.66"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("E24").Value = "  -7.04%  "
$ws.Range("E25").Value = "  -8.08%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.21"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.04"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.13"
$ws.Range("E29").Value = "  -3.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.23"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.22"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.97"
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.58"
$ws.Range("E33").Value = "  +11.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0832"
$ws.Range("E34").Value = "  -2.62%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -3.95%  "
$ws.Range("E38").Value = "  -3.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.94"
$ws.Range("E39").Value = "  +4.16%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.92"
$ws.Range("E40").Value = "  -8.34%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.42"
$ws.Range("E41").Value = "  -10.24%  "
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.777.61"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "87.50"
$ws.Range("E45").Value = "  +5.62%  "
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.06"
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.32"
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.58"
$ws.Range("E50").Value = "  -5.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.80"
$ws.Range("E51").Value = "  -5.61%  "
